$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data rows down by inserting a new (blank) row right
# below the header row, then retarget the "/" endpoint's Result text
# onto the row that slides into that slot.
$ws.Rows.Item(2).Insert()

$ws.Range("A3").Value = "/"
$ws.Range("C3").Value = "logout page"

# Append a new row at the bottom of the table for the "/login" endpoint.
$ws.Range("A9").Value = "/login"
$ws.Range("B9").Value = "GET, POST"
$ws.Range("C9").Value = "Ввод логина и пароля"

# Update the active cell selection to match the authored state.
$ws.Range("B7").Select()
